$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 2
$ws.Range("F2").Value = 81
$ws.Range("H2").Value = 'kitchens'
$ws.Range("L2").Value = 'stimuli/img_ifebc.png'
$ws.Range("M2").Value = 84
$ws.Range("N2").Value = 65.88235294117646
$ws.Range("O2").Value = 74.94117647058823
$ws.Range("P2").Value = 34
$ws.Range("R2").Value = 9
$ws.Range("S2").Value = 9
$ws.Range("C3").Value = 2
$ws.Range("F3").Value = 82
$ws.Range("H3").Value = 'kitchens'
$ws.Range("I3").Value = $null
$ws.Range("J3").Value = 'new'
$ws.Range("K3").Value = 'f'
$ws.Range("L3").Value = 'stimuli/img_mgnmm.png'
$ws.Range("M3").Value = 79.1470588235294
$ws.Range("N3").Value = 60.38235294117647
$ws.Range("O3").Value = 69.76470588235294
$ws.Range("P3").Value = 34
$ws.Range("Q3").Value = 8
$ws.Range("R3").Value = 8
$ws.Range("S3").Value = 8
$ws.Range("C4").Value = 2
$ws.Range("F4").Value = 83
$ws.Range("H4").Value = 'kitchens'
$ws.Range("I4").Value = 'target'
$ws.Range("J4").Value = 'old'
$ws.Range("K4").Value = 'j'
$ws.Range("L4").Value = 'stimuli/img_c0me7.png'
$ws.Range("M4").Value = 68.4
$ws.Range("N4").Value = 45.62857142857143
$ws.Range("O4").Value = 57.01428571428572
$ws.Range("P4").Value = 35
$ws.Range("Q4").Value = 4
$ws.Range("R4").Value = 4
$ws.Range("S4").Value = 4
$ws.Range("C5").Value = 2
$ws.Range("F5").Value = 84
$ws.Range("H5").Value = 'kitchens'
$ws.Range("L5").Value = 'stimuli/img_es7o2.png'
$ws.Range("M5").Value = 52.48571428571429
$ws.Range("N5").Value = 27.54285714285714
$ws.Range("O5").Value = 40.01428571428572
$ws.Range("P5").Value = 35
$ws.Range("Q5").Value = 2
$ws.Range("R5").Value = 2
$ws.Range("S5").Value = 2
$ws.Range("C6").Value = 2
$ws.Range("F6").Value = 85
$ws.Range("H6").Value = 'kitchens'
$ws.Range("L6").Value = 'stimuli/img_05flq.png'
$ws.Range("M6").Value = 47.10344827586207
$ws.Range("N6").Value = 25.72413793103448
$ws.Range("O6").Value = 36.41379310344828
$ws.Range("P6").Value = 29
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = 1
$ws.Range("S6").Value = 1
$ws.Range("C7").Value = 2
$ws.Range("F7").Value = 86
$ws.Range("H7").Value = 'kitchens'
$ws.Range("L7").Value = 'stimuli/img_p659z.png'
$ws.Range("M7").Value = 84.21621621621621
$ws.Range("N7").Value = 65.37837837837837
$ws.Range("O7").Value = 74.79729729729729
$ws.Range("P7").Value = 37
$ws.Range("Q7").Value = 9
$ws.Range("R7").Value = 9
$ws.Range("S7").Value = 9
$ws.Range("C8").Value = 2
$ws.Range("F8").Value = 87
$ws.Range("H8").Value = 'kitchens'
$ws.Range("L8").Value = 'stimuli/img_xesl0.png'
$ws.Range("M8").Value = 69.28571428571429
$ws.Range("N8").Value = 47.35714285714285
$ws.Range("O8").Value = 58.32142857142857
$ws.Range("P8").Value = 28
$ws.Range("Q8").Value = 5
$ws.Range("R8").Value = 5
$ws.Range("S8").Value = 5
$ws.Range("C9").Value = 2
$ws.Range("F9").Value = 88
$ws.Range("H9").Value = 'kitchens'
$ws.Range("L9").Value = 'stimuli/img_bwo9g.png'
$ws.Range("M9").Value = 64.81818181818181
$ws.Range("N9").Value = 42.36363636363637
$ws.Range("O9").Value = 53.59090909090909
$ws.Range("P9").Value = 33
$ws.Range("Q9").Value = 4
$ws.Range("R9").Value = 4
$ws.Range("S9").Value = 4
$ws.Range("C10").Value = 2
$ws.Range("F10").Value = 89
$ws.Range("H10").Value = 'kitchens'
$ws.Range("L10").Value = 'stimuli/img_z5osu.png'
$ws.Range("M10").Value = 71.42857142857143
$ws.Range("N10").Value = 47.34285714285714
$ws.Range("O10").Value = 59.38571428571429
$ws.Range("P10").Value = 35
$ws.Range("Q10").Value = 5
$ws.Range("R10").Value = 5
$ws.Range("S10").Value = 5
$ws.Range("C11").Value = 2
$ws.Range("F11").Value = 90
$ws.Range("H11").Value = 'kitchens'
$ws.Range("L11").Value = 'stimuli/img_zi8qc.png'
$ws.Range("M11").Value = 77.14285714285714
$ws.Range("N11").Value = 57.02857142857143
$ws.Range("O11").Value = 67.08571428571429
$ws.Range("P11").Value = 35
$ws.Range("Q11").Value = 7
$ws.Range("R11").Value = 7
$ws.Range("S11").Value = 7
$ws.Range("C12").Value = 2
$ws.Range("F12").Value = 91
$ws.Range("H12").Value = 'kitchens'
$ws.Range("I12").Value = $null
$ws.Range("J12").Value = 'new'
$ws.Range("K12").Value = 'f'
$ws.Range("L12").Value = 'stimuli/img_mawe6.png'
$ws.Range("M12").Value = 83.48387096774194
$ws.Range("N12").Value = 65.54838709677419
$ws.Range("O12").Value = 74.51612903225806
$ws.Range("P12").Value = 31
$ws.Range("Q12").Value = 9
$ws.Range("R12").Value = 9
$ws.Range("S12").Value = 9
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 92
$ws.Range("H13").Value = 'kitchens'
$ws.Range("I13").Value = $null
$ws.Range("J13").Value = 'new'
$ws.Range("K13").Value = 'f'
$ws.Range("L13").Value = 'stimuli/img_vbrb7.png'
$ws.Range("M13").Value = 85.5625
$ws.Range("N13").Value = 71.46875
$ws.Range("O13").Value = 78.515625
$ws.Range("P13").Value = 32
$ws.Range("Q13").Value = 10
$ws.Range("R13").Value = 10
$ws.Range("S13").Value = 10
$ws.Range("C14").Value = 2
$ws.Range("F14").Value = 93
$ws.Range("H14").Value = 'kitchens'
$ws.Range("L14").Value = 'stimuli/img_ua9bs.png'
$ws.Range("M14").Value = 82
$ws.Range("N14").Value = 62.23333333333333
$ws.Range("O14").Value = 72.11666666666667
$ws.Range("P14").Value = 30
$ws.Range("Q14").Value = 9
$ws.Range("R14").Value = 9
$ws.Range("S14").Value = 9
$ws.Range("C15").Value = 2
$ws.Range("F15").Value = 94
$ws.Range("H15").Value = 'kitchens'
$ws.Range("L15").Value = 'stimuli/img_i2k07.png'
$ws.Range("M15").Value = 64.25925925925925
$ws.Range("N15").Value = 40.92592592592592
$ws.Range("O15").Value = 52.59259259259259
$ws.Range("P15").Value = 27
$ws.Range("Q15").Value = 3
$ws.Range("R15").Value = 3
$ws.Range("S15").Value = 3
$ws.Range("C16").Value = 2
$ws.Range("F16").Value = 95
$ws.Range("H16").Value = 'kitchens'
$ws.Range("L16").Value = 'stimuli/img_ncr40.png'
$ws.Range("M16").Value = 75.66666666666667
$ws.Range("N16").Value = 54.27272727272727
$ws.Range("O16").Value = 64.96969696969697
$ws.Range("Q16").Value = 6
$ws.Range("R16").Value = 6
$ws.Range("S16").Value = 6
$ws.Range("C17").Value = 2
$ws.Range("F17").Value = 96
$ws.Range("H17").Value = 'kitchens'
$ws.Range("I17").Value = $null
$ws.Range("J17").Value = 'new'
$ws.Range("K17").Value = 'f'
$ws.Range("L17").Value = 'stimuli/img_qmgwq.png'
$ws.Range("M17").Value = 84.58333333333333
$ws.Range("N17").Value = 64.44444444444444
$ws.Range("O17").Value = 74.51388888888889
$ws.Range("P17").Value = 36
$ws.Range("Q17").Value = 9
$ws.Range("R17").Value = 9
$ws.Range("S17").Value = 9
$ws.Range("C18").Value = 2
$ws.Range("F18").Value = 97
$ws.Range("H18").Value = $null
$ws.Range("I18").Value = $null
$ws.Range("J18").Value = 'catch'
$ws.Range("K18").Value = 'f'
$ws.Range("L18").Value = 'stimuli/catch_20.jpg'
$ws.Range("M18").Value = $null
$ws.Range("N18").Value = $null
$ws.Range("O18").Value = $null
$ws.Range("P18").Value = $null
$ws.Range("Q18").Value = $null
$ws.Range("R18").Value = $null
$ws.Range("S18").Value = $null
$ws.Range("C19").Value = 2
$ws.Range("F19").Value = 98
$ws.Range("H19").Value = 'kitchens'
$ws.Range("I19").Value = 'target'
$ws.Range("J19").Value = 'old'
$ws.Range("K19").Value = 'j'
$ws.Range("L19").Value = 'stimuli/img_jz3kd.png'
$ws.Range("M19").Value = 72.79411764705883
$ws.Range("N19").Value = 51.64705882352941
$ws.Range("O19").Value = 62.22058823529412
$ws.Range("P19").Value = 34
$ws.Range("Q19").Value = 6
$ws.Range("R19").Value = 6
$ws.Range("S19").Value = 6
$ws.Range("C20").Value = 2
$ws.Range("F20").Value = 99
$ws.Range("H20").Value = 'kitchens'
$ws.Range("L20").Value = 'stimuli/img_cv9qj.png'
$ws.Range("M20").Value = 60.34375
$ws.Range("N20").Value = 35.34375
$ws.Range("O20").Value = 47.84375
$ws.Range("P20").Value = 32
$ws.Range("Q20").Value = 3
$ws.Range("R20").Value = 3
$ws.Range("S20").Value = 3
$ws.Range("C21").Value = 2
$ws.Range("F21").Value = 100
$ws.Range("H21").Value = 'kitchens'
$ws.Range("L21").Value = 'stimuli/img_uy1n4.png'
$ws.Range("M21").Value = 76.30555555555556
$ws.Range("N21").Value = 55.33333333333334
$ws.Range("O21").Value = 65.81944444444444
$ws.Range("P21").Value = 36
$ws.Range("Q21").Value = 7
$ws.Range("R21").Value = 7
$ws.Range("S21").Value = 7
$ws.Range("C22").Value = 2
$ws.Range("F22").Value = 101
$ws.Range("H22").Value = 'kitchens'
$ws.Range("L22").Value = 'stimuli/img_c79r7.png'
$ws.Range("M22").Value = 56.26470588235294
$ws.Range("N22").Value = 34.26470588235294
$ws.Range("O22").Value = 45.26470588235294
$ws.Range("P22").Value = 34
$ws.Range("Q22").Value = 2
$ws.Range("R22").Value = 2
$ws.Range("S22").Value = 2
$ws.Range("C23").Value = 2
$ws.Range("F23").Value = 102
$ws.Range("H23").Value = 'kitchens'
$ws.Range("L23").Value = 'stimuli/img_84s7n.png'
$ws.Range("M23").Value = 11.03125
$ws.Range("N23").Value = 2.90625
$ws.Range("O23").Value = 6.96875
$ws.Range("P23").Value = 32
$ws.Range("Q23").Value = 1
$ws.Range("R23").Value = 1
$ws.Range("S23").Value = 1
$ws.Range("C24").Value = 2
$ws.Range("F24").Value = 103
$ws.Range("H24").Value = 'kitchens'
$ws.Range("L24").Value = 'stimuli/img_xti0z.png'
$ws.Range("M24").Value = 81.40625
$ws.Range("N24").Value = 61.4375
$ws.Range("O24").Value = 71.421875
$ws.Range("P24").Value = 32
$ws.Range("Q24").Value = 8
$ws.Range("R24").Value = 8
$ws.Range("S24").Value = 8
$ws.Range("C25").Value = 2
$ws.Range("F25").Value = 104
$ws.Range("H25").Value = 'kitchens'
$ws.Range("I25").Value = 'target'
$ws.Range("J25").Value = 'old'
$ws.Range("K25").Value = 'j'
$ws.Range("L25").Value = 'stimuli/img_7pgd2.png'
$ws.Range("M25").Value = 78.59375
$ws.Range("N25").Value = 57.84375
$ws.Range("O25").Value = 68.21875
$ws.Range("P25").Value = 32
$ws.Range("Q25").Value = 8
$ws.Range("R25").Value = 7
$ws.Range("S25").Value = 7
$ws.Range("C26").Value = 2
$ws.Range("F26").Value = 105
$ws.Range("H26").Value = 'kitchens'
$ws.Range("L26").Value = 'stimuli/img_s9are.png'
$ws.Range("M26").Value = 90.14285714285714
$ws.Range("N26").Value = 75.22857142857143
$ws.Range("O26").Value = 82.68571428571428
$ws.Range("P26").Value = 35
$ws.Range("Q26").Value = 10
$ws.Range("R26").Value = 10
$ws.Range("S26").Value = 10
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 106
$ws.Range("H27").Value = 'kitchens'
$ws.Range("I27").Value = 'target'
$ws.Range("J27").Value = 'old'
$ws.Range("K27").Value = 'j'
$ws.Range("L27").Value = 'stimuli/img_6nbgt.png'
$ws.Range("M27").Value = 78.45161290322581
$ws.Range("N27").Value = 57.83870967741935
$ws.Range("O27").Value = 68.14516129032258
$ws.Range("P27").Value = 31
$ws.Range("Q27").Value = 7
$ws.Range("R27").Value = 7
$ws.Range("S27").Value = 7
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 107
$ws.Range("H28").Value = 'kitchens'
$ws.Range("L28").Value = 'stimuli/img_hfz8w.png'
$ws.Range("M28").Value = 55.46153846153846
$ws.Range("N28").Value = 27.28205128205128
$ws.Range("O28").Value = 41.37179487179487
$ws.Range("P28").Value = 39
$ws.Range("Q28").Value = 2
$ws.Range("R28").Value = 2
$ws.Range("S28").Value = 2
$ws.Range("C29").Value = 2
$ws.Range("F29").Value = 108
$ws.Range("H29").Value = 'kitchens'
$ws.Range("I29").Value = 'target'
$ws.Range("J29").Value = 'old'
$ws.Range("K29").Value = 'j'
$ws.Range("L29").Value = 'stimuli/img_j5rpx.png'
$ws.Range("M29").Value = 72.24242424242425
$ws.Range("N29").Value = 50
$ws.Range("O29").Value = 61.12121212121212
$ws.Range("P29").Value = 33
$ws.Range("Q29").Value = 5
$ws.Range("R29").Value = 5
$ws.Range("S29").Value = 5
